$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.296.54'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.00%  '
$ws.Range("E2").NumberFormat = "General"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.677.63'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.71%  '
$ws.Range("E3").NumberFormat = "General"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("E4").NumberFormat = "General"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.73'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.48%  '
$ws.Range("E5").NumberFormat = "General"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5256'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.16%  '
$ws.Range("E6").NumberFormat = "General"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.38%  '
$ws.Range("E7").NumberFormat = "General"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.11%  '
$ws.Range("E8").NumberFormat = "General"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06458'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.96%  '
$ws.Range("E9").NumberFormat = "General"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.83'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.45%  '
$ws.Range("E10").NumberFormat = "General"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07500'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.23%  '
$ws.Range("E11").NumberFormat = "General"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.694.60'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.68%  '
$ws.Range("E12").NumberFormat = "General"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.512'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.24%  '
$ws.Range("E13").NumberFormat = "General"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5767'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.78%  '
$ws.Range("E14").NumberFormat = "General"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008483'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.52%  '
$ws.Range("E15").NumberFormat = "General"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.63'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.58%  '
$ws.Range("E16").NumberFormat = "General"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.328.48'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.90%  '
$ws.Range("E17").NumberFormat = "General"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.907'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.10%  '
$ws.Range("E18").NumberFormat = "General"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.87'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.64%  '
$ws.Range("E20").NumberFormat = "General"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '189.41'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.28%  '
$ws.Range("E21").NumberFormat = "General"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.174'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.45%  '
$ws.Range("E22").NumberFormat = "General"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.31%  '
$ws.Range("E23").NumberFormat = "General"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '144.97'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.51%  '
$ws.Range("E24").NumberFormat = "General"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.22%  '
$ws.Range("E25").NumberFormat = "General"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1257'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +5.72%  '
$ws.Range("E26").NumberFormat = "General"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.78%  '
$ws.Range("E27").NumberFormat = "General"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06439'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.48%  '
$ws.Range("E28").NumberFormat = "General"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.51%  '
$ws.Range("E29").NumberFormat = "General"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.75%  '
$ws.Range("E30").NumberFormat = "General"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.580'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.59%  '
$ws.Range("E31").NumberFormat = "General"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.583'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.30%  '
$ws.Range("E32").NumberFormat = "General"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.70%  '
$ws.Range("E33").NumberFormat = "General"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.026'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.78%  '
$ws.Range("E34").NumberFormat = "General"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6186'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.11%  '
$ws.Range("E35").NumberFormat = "General"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.407'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.74%  '
$ws.Range("E36").NumberFormat = "General"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.739'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.10%  '
$ws.Range("E37").NumberFormat = "General"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.283'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.13%  '
$ws.Range("E38").NumberFormat = "General"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.114.79'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.68%  '
$ws.Range("E39").NumberFormat = "General"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.62%  '
$ws.Range("E40").NumberFormat = "General"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8713'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.35%  '
$ws.Range("E41").NumberFormat = "General"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.71%  '
$ws.Range("E42").NumberFormat = "General"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("E43").NumberFormat = "General"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.828.40'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.85%  '
$ws.Range("E44").NumberFormat = "General"

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Aave'
$ws.Range("B45").NumberFormat = "General"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("C45").NumberFormat = "General"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '56.87'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.11%  '
$ws.Range("E45").NumberFormat = "General"

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("B46").NumberFormat = "General"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("C46").NumberFormat = "General"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000106'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -8.41%  '
$ws.Range("E46").NumberFormat = "General"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.008'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.17%  '
$ws.Range("E47").NumberFormat = "General"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.132'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.36%  '
$ws.Range("E48").NumberFormat = "General"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05264'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.00%  '
$ws.Range("E49").NumberFormat = "General"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.13%  '
$ws.Range("E50").NumberFormat = "General"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.043'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.50%  '
$ws.Range("E51").NumberFormat = "General"
